# Apply the day-to-day code coverage update:
# - add the "total" column (E) percentages for rows 2-6,8
# - update window size in workbook view
# - update active selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column E values (percentage formatted like column D, style s="2")
$ws.Range("E2").Value = 0.701
$ws.Range("E3").Value = 0.034
$ws.Range("E4").Value = 0.034
$ws.Range("E5").Value = 0.171
$ws.Range("E6").Value = 0.379
$ws.Range("E8").Value = 0.177

# Match the number format used by column D (percentage with two decimals, numFmtId 10)
$ws.Range("E2:E6").NumberFormat = "0.00%"
$ws.Range("E8").NumberFormat = "0.00%"

# Size the new column to fit its content, like the existing bestFit columns
$ws.Columns.Item(5).ColumnWidth = 8.3

# Update the selection to G7 as in the final sheet
$ws.Range("G7").Select()

# Update the workbook window size
$excel.Width = 24180
$excel.Height = 4980
